$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns remain text so numeric-looking values
# (e.g. "1.006", "0.00001080") are not silently coerced to numbers/
# scientific notation by the COM layer.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '22.169.47'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '1.556.19'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").Value = '286.36'
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").Value = '0.3803'
$ws.Range("E7").Value = '  +3.86%  '
$ws.Range("D8").Value = '0.3265'
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").Value = '43.92'
$ws.Range("E9").Value = '  -9.49%  '
$ws.Range("D10").Value = '1.133'
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '0.07395'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '1.005'
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").Value = '20.38'
$ws.Range("E13").Value = '  -2.43%  '
$ws.Range("D14").Value = '5.823'
$ws.Range("E14").Value = '  -2.94%  '
$ws.Range("D15").Value = '6.782'
$ws.Range("E15").Value = '  -2.17%  '
$ws.Range("D16").Value = '1.559.56'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '0.00001080'
$ws.Range("E17").Value = '  -2.80%  '
$ws.Range("D18").Value = '0.06692'
$ws.Range("E18").Value = '  -1.02%  '
$ws.Range("D19").Value = '85.84'
$ws.Range("E19").Value = '  -2.83%  '
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '6.372'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("E22").Value = '  -1.97%  '
$ws.Range("D23").Value = '11.63'
$ws.Range("E23").Value = '  -4.33%  '
$ws.Range("D24").Value = '22.350.16'
$ws.Range("E24").Value = '  -0.62%  '
$ws.Range("D25").Value = '2.296'
$ws.Range("E25").Value = '  -3.99%  '
$ws.Range("D26").Value = '2.566'
$ws.Range("E26").Value = '  -2.18%  '
$ws.Range("D27").Value = '149.99'
$ws.Range("E27").Value = '  -1.81%  '
$ws.Range("D28").Value = '19.52'
$ws.Range("E28").Value = '  -0.54%  '
$ws.Range("D29").Value = '4.920'
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").Value = '122.48'
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("D31").Value = '1.741.38'
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").Value = '1.067'
$ws.Range("E32").Value = '  +1.75%  '
$ws.Range("D33").Value = '5.968'
$ws.Range("E33").Value = '  -3.42%  '
$ws.Range("D34").Value = '1.859'
$ws.Range("E34").Value = '  -7.20%  '
$ws.Range("D35").Value = '9.489'
$ws.Range("E35").Value = '  -3.73%  '
$ws.Range("D36").Value = '0.08250'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("D37").Value = '0.02368'
$ws.Range("E37").Value = '  -3.21%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.278'
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '5.285'
$ws.Range("E39").Value = '  -2.94%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = '0.06233'
$ws.Range("E40").Value = '  -3.67%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2165'
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("D42").Value = '11.05'
$ws.Range("E42").Value = '  -2.90%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.6088'
$ws.Range("E43").Value = '  -4.15%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '13.67'
$ws.Range("E45").Value = '  -2.25%  '
$ws.Range("D46").Value = '3.745'
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").Value = '0.5903'
$ws.Range("E47").Value = '  -4.11%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '1.994'
$ws.Range("E48").Value = '  -3.09%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '123.22'
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("D50").Value = '1.171'
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("D51").Value = '0.07083'
$ws.Range("E51").Value = '  -2.23%  '
